# ore_tirocinio.docx edit script
# - merges the " parte 4" + "." runs (cosmetic, no text change)
# - inserts " banda colori" into the "...e modifica del raster..." sentence
# - removes "prodotto " from "...raster prodotto su QGIS."
# - appends a new activity row (12/10/2023) plus a trailing blank row
# - bumps the "Totale ore" total from 46 to 52

$d = $word.ActiveDocument

# Useful special characters
$enDash   = [char]0x2013
$lq       = [char]0x201C
$rq       = [char]0x201D

# --- 1) Merge " parte 4" + "." into a single run (text itself is unchanged) ---
$d.Content.Find.Execute(" parte 4.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " parte 4.", 2) | Out-Null

# --- 2) Insert " banda colori" before " del raster" ---
$d.Content.Find.Execute("e modifica del ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "e modifica banda colori del ", 2) | Out-Null

# --- 3) Drop "prodotto " so the sentence reads "...raster su QGIS." ---
$d.Content.Find.Execute(" prodotto su QGIS.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " su QGIS.", 2) | Out-Null

# --- 4) Add the two new table rows (content row + trailing blank row) ---
$tbl = $d.Tables.Item(1)
$lastRow = $tbl.Rows.Item($tbl.Rows.Count)
$contentRow = $tbl.Rows.Add($lastRow)
$blankRow = $tbl.Rows.Add($lastRow)

$newRowIndex = $contentRow.Index

$tbl.Cell($newRowIndex, 1).Range.Text = "12/10/2023"
$tbl.Cell($newRowIndex, 2).Range.Text = "9:30 " + $enDash + " 15:30"
$tbl.Cell($newRowIndex, 3).Range.Text = "6"
$tbl.Cell($newRowIndex, 4).Range.Text = "Creazione dello script " + $lq + "generatore_grafico_ndvi_campania" + $rq + ".js"

# --- 5) Update the hour total ---
$d.Content.Find.Execute("Totale ore: 46", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Totale ore: 52", 2) | Out-Null
